# Update the lattice multiplication exercise table: replace the contents
# of every cell with new "A x B" problems (and their corresponding lattice
# scaffolding lines), row by row, left to right.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$nl = [char]11  # line-break within a run (maps to <w:br/> between <w:t> runs)

$cellUpdates = @(
    @{ Row = 1; Col = 1; Text = ("79 x 43" + $nl + "  4    3" + $nl + "  ----" + $nl + "7|    |" + $nl + "9|    |") },
    @{ Row = 1; Col = 2; Text = ("26 x 18" + $nl + "  1    8" + $nl + "  ----" + $nl + "2|    |" + $nl + "6|    |") },
    @{ Row = 1; Col = 3; Text = ("76 x 53" + $nl + "  5    3" + $nl + "  ----" + $nl + "7|    |" + $nl + "6|    |") },
    @{ Row = 2; Col = 1; Text = ("55 x 41" + $nl + "  4    1" + $nl + "  ----" + $nl + "5|    |" + $nl + "5|    |") },
    @{ Row = 2; Col = 2; Text = ("94 x 58" + $nl + "  5    8" + $nl + "  ----" + $nl + "9|    |" + $nl + "4|    |") },
    @{ Row = 2; Col = 3; Text = ("65 x 40" + $nl + "  4    0" + $nl + "  ----" + $nl + "6|    |" + $nl + "5|    |") },
    @{ Row = 3; Col = 1; Text = ("21 x 53" + $nl + "  5    3" + $nl + "  ----" + $nl + "2|    |" + $nl + "1|    |") },
    @{ Row = 3; Col = 2; Text = ("87 x 67" + $nl + "  6    7" + $nl + "  ----" + $nl + "8|    |" + $nl + "7|    |") },
    @{ Row = 3; Col = 3; Text = ("59 x 13" + $nl + "  1    3" + $nl + "  ----" + $nl + "5|    |" + $nl + "9|    |") },
    @{ Row = 4; Col = 1; Text = ("62 x 58" + $nl + "  5    8" + $nl + "  ----" + $nl + "6|    |" + $nl + "2|    |") },
    @{ Row = 4; Col = 2; Text = ("80 x 36" + $nl + "  3    6" + $nl + "  ----" + $nl + "8|    |" + $nl + "0|    |") },
    @{ Row = 4; Col = 3; Text = ("12 x 35" + $nl + "  3    5" + $nl + "  ----" + $nl + "1|    |" + $nl + "2|    |") },
    @{ Row = 5; Col = 1; Text = ("59 x 98" + $nl + "  9    8" + $nl + "  ----" + $nl + "5|    |" + $nl + "9|    |") },
    @{ Row = 5; Col = 2; Text = ("34 x 89" + $nl + "  8    9" + $nl + "  ----" + $nl + "3|    |" + $nl + "4|    |") },
    @{ Row = 5; Col = 3; Text = ("20 x 20" + $nl + "  2    0" + $nl + "  ----" + $nl + "2|    |" + $nl + "0|    |") }
)

foreach ($u in $cellUpdates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}

Write-Host "Updated" $cellUpdates.Count "cells"
